$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'symptom_group'
$ws.Range("B1").Value = 'ABG+UNKNOWN'
$ws.Range("C1").Value = 'ABG+VBG'
$ws.Range("D1").Value = 'ABG+VBG+UNKNOWN'
$ws.Range("E1").Value = 'ABG-only'
$ws.Range("F1").Value = 'UNKNOWN-only'
$ws.Range("G1").Value = 'VBG+UNKNOWN'
$ws.Range("H1").Value = 'VBG-only'

$ws.Range("A2").Value = 'Administrative'
$ws.Range("B2").Value = 3.4
$ws.Range("C2").Value = 1.8
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 2.4
$ws.Range("F2").Value = 1.9
$ws.Range("G2").Value = 3.1
$ws.Range("H2").Value = 1.4

$ws.Range("A3").Value = 'Diseases (patient-stated)'
$ws.Range("B3").Value = 3.4
$ws.Range("C3").Value = 4.9
$ws.Range("D3").Value = 7.1
$ws.Range("E3").Value = 7.7
$ws.Range("F3").Value = 5.7
$ws.Range("G3").Value = 7.8
$ws.Range("H3").Value = 4.1

$ws.Range("A4").Value = 'Injuries & adverse effects'
$ws.Range("B4").Value = 24.6
$ws.Range("C4").Value = 16.6
$ws.Range("D4").Value = 14.3
$ws.Range("E4").Value = 22.7
$ws.Range("F4").Value = 15.2
$ws.Range("G4").Value = 6.2
$ws.Range("H4").Value = 10.4

$ws.Range("A5").Value = 'Other'
$ws.Range("B5").Value = 5.9
$ws.Range("C5").Value = 3.3
$ws.Range("D5").Value = 10.7
$ws.Range("E5").Value = 6.3
$ws.Range("F5").Value = 7.6
$ws.Range("G5").Value = 4.7
$ws.Range("H5").Value = 6.6

$ws.Range("A6").Value = 'Symptom – Circulatory'
$ws.Range("B6").Value = 9.3
$ws.Range("C6").Value = 6.3
$ws.Range("D6").Value = 7.1
$ws.Range("E6").Value = 6.4
$ws.Range("F6").Value = 8.9
$ws.Range("G6").Value = 10.9
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 'Symptom – Digestive'
$ws.Range("B7").Value = 11.9
$ws.Range("C7").Value = 7.3
$ws.Range("D7").Value = 7.1
$ws.Range("E7").Value = 10.6
$ws.Range("F7").Value = 12.7
$ws.Range("G7").Value = 3.1
$ws.Range("H7").Value = 10.1

$ws.Range("A8").Value = 'Symptom – General'
$ws.Range("B8").Value = 2.5
$ws.Range("C8").Value = 2.4
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 4.4
$ws.Range("G8").Value = 7.8
$ws.Range("H8").Value = 4.1

$ws.Range("A9").Value = 'Symptom – Nervous'
$ws.Range("B9").Value = 5.9
$ws.Range("C9").Value = 12.3
$ws.Range("D9").Value = 10.7
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 16.5
$ws.Range("G9").Value = 10.9
$ws.Range("H9").Value = 13.5

$ws.Range("A10").Value = 'Symptom – Respiratory'
$ws.Range("B10").Value = 24.6
$ws.Range("C10").Value = 41.6
$ws.Range("D10").Value = 41.1
$ws.Range("E10").Value = 24
$ws.Range("F10").Value = 23.4
$ws.Range("G10").Value = 43.8
$ws.Range("H10").Value = 38.4

$ws.Range("A11").Value = 'Symptom – Skin/Hair/Nails'
$ws.Range("B11").Value = 3.4
$ws.Range("C11").Value = 0.8
$ws.Range("D11").Value = 1.8
$ws.Range("E11").Value = 2.7
$ws.Range("F11").Value = 2.5
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1

$ws.Range("A12").Value = 'Uncodable/Unknown'
$ws.Range("B12").Value = 5.1
$ws.Range("C12").Value = 2.6
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 4.4
$ws.Range("F12").Value = 1.3
$ws.Range("G12").Value = 1.6
$ws.Range("H12").Value = 1.3

